# Apply the edit described in the diff:
#  - Rows 42-72, column C: "Propia (CCNN)" -> "Propia (Empalme CCNN)"
#  - Rows 49-70, column B: updated numeric cost values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (author label) for rows 42 through 72
for ($row = 42; $row -le 72; $row++) {
    $ws.Cells.Item($row, 3).Value = "Propia (Empalme CCNN)"
}

# Update column B (cost values) for rows 49 through 70
$ws.Cells.Item(49, 2).Value = 10.02154264552858
$ws.Cells.Item(50, 2).Value = 6.564068037090995
$ws.Cells.Item(51, 2).Value = 9.156910604924347
$ws.Cells.Item(52, 2).Value = 14.03924140416419
$ws.Cells.Item(53, 2).Value = 11.18015821598729
$ws.Cells.Item(54, 2).Value = 6.039871827719271
$ws.Cells.Item(55, 2).Value = 7.466465050297496
$ws.Cells.Item(56, 2).Value = 7.636867481217681
$ws.Cells.Item(57, 2).Value = 8.803083440442739
$ws.Cells.Item(58, 2).Value = 11.65414131088522
$ws.Cells.Item(59, 2).Value = 14.09724983579705
$ws.Cells.Item(60, 2).Value = 16.87092470750025
$ws.Cells.Item(61, 2).Value = 14.82161676785997
$ws.Cells.Item(62, 2).Value = 17.62837777096416
$ws.Cells.Item(63, 2).Value = 23.0656320320991
$ws.Cells.Item(64, 2).Value = 26.38578606501133
$ws.Cells.Item(65, 2).Value = 29.58253788187292
$ws.Cells.Item(66, 2).Value = 25.59702003976075
$ws.Cells.Item(67, 2).Value = 24.26300522059164
$ws.Cells.Item(68, 2).Value = 21.97031226820259
$ws.Cells.Item(69, 2).Value = 23.10442440072572
$ws.Cells.Item(70, 2).Value = 17.54025290647301
